# Update cryptos list prices/volume changes (GitHub Actions scheduled refresh).
# Sets cell values as text (preserving formatting such as thousands separators,
# padded percentage strings, and the subscript-zero notation) by temporarily
# forcing a text number format, assigning the value, then restoring the
# cell's original style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "28.495.95"
Set-TextValue $ws "E2" "  +4.74%  "
Set-TextValue $ws "D3" "1.590.89"
Set-TextValue $ws "E4" "  -0.74%  "
Set-TextValue $ws "D5" "214.82"
Set-TextValue $ws "E5" "  +1.75%  "
Set-TextValue $ws "E6" "  +1.08%  "
Set-TextValue $ws "E7" "  -0.66%  "
Set-TextValue $ws "D8" "23.84"
Set-TextValue $ws "E8" "  +8.09%  "
Set-TextValue $ws "E9" "  +1.44%  "
Set-TextValue $ws "E10" "  +0.39%  "
Set-TextValue $ws "E11" "  +2.21%  "
Set-TextValue $ws "D12" "1.817.88"
Set-TextValue $ws "E12" "  +1.36%  "
Set-TextValue $ws "D13" "1.613.88"
Set-TextValue $ws "E13" "  +2.74%  "
Set-TextValue $ws "E14" "  +0.14%  "
Set-TextValue $ws "E15" "  +2.09%  "
Set-TextValue $ws "D16" "28.486.01"
Set-TextValue $ws "E16" "  +4.93%  "
Set-TextValue $ws "D17" "64.11"
Set-TextValue $ws "E17" "  +3.07%  "
Set-TextValue $ws "D18" "233.16"
Set-TextValue $ws "E18" "  +7.61%  "
Set-TextValue $ws "D19" "7.53"
Set-TextValue $ws "E19" "  +0.41%  "
Set-TextValue $ws "E20" "  +0.99%  "
Set-TextValue $ws "E21" "  -0.65%  "
Set-TextValue $ws "E22" "  -0.34%  "
Set-TextValue $ws "D23" "9.41"
Set-TextValue $ws "E23" "  +1.96%  "
Set-TextValue $ws "D24" "1.95"
Set-TextValue $ws "E24" "  +0.28%  "
Set-TextValue $ws "D25" "152.25"
Set-TextValue $ws "E25" "  -1.12%  "
Set-TextValue $ws "D26" "15.33"
Set-TextValue $ws "E26" "  +1.61%  "
Set-TextValue $ws "D27" "6.63"
Set-TextValue $ws "E27" "  -0.37%  "
Set-TextValue $ws "D28" "0.108"
Set-TextValue $ws "E28" "  +1.03%  "
Set-TextValue $ws "E29" "  -0.66%  "
Set-TextValue $ws "E30" "  +0.38%  "
Set-TextValue $ws "E31" "  +0.35%  "
Set-TextValue $ws "D32" "3.25"
Set-TextValue $ws "E32" "  +0.00%  "
Set-TextValue $ws "E33" "  -0.28%  "
Set-TextValue $ws "D34" "1.421.62"
Set-TextValue $ws "E34" "  -1.93%  "
Set-TextValue $ws "E35" "  -0.23%  "
Set-TextValue $ws "E36" "  -6.12%  "
Set-TextValue $ws "E37" "  -0.58%  "
Set-TextValue $ws "E38" "  +0.32%  "
Set-TextValue $ws "D39" "2.52"
Set-TextValue $ws "E39" "  +7.79%  "
Set-TextValue $ws "D40" "0.544"
Set-TextValue $ws "E40" "  +2.11%  "
Set-TextValue $ws "D41" "0.815"
Set-TextValue $ws "E41" "  +0.56%  "
Set-TextValue $ws "E42" "  -0.77%  "
Set-TextValue $ws "D43" "5.67"
Set-TextValue $ws "E43" "  -2.91%  "
Set-TextValue $ws "B44" "RenderToken"
Set-TextValue $ws "C44" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D44" "1.83"
Set-TextValue $ws "E44" "  +6.00%  "
Set-TextValue $ws "B45" "WEMIXToken"
Set-TextValue $ws "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D45" "0.977"
Set-TextValue $ws "E45" "  -2.75%  "
Set-TextValue $ws "D46" "64.70"
Set-TextValue $ws "E46" "  +0.42%  "
Set-TextValue $ws "D47" "1.727.54"
Set-TextValue $ws "E47" "  +1.29%  "
Set-TextValue $ws "D48" "87.86"
Set-TextValue $ws "E48" "  +2.07%  "
Set-TextValue $ws "B49" "BabyDogeCoin"
Set-TextValue $ws "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D49" "0.0₆0107"
Set-TextValue $ws "E49" "  +2.37%  "
Set-TextValue $ws "B50" "Cronos"
Set-TextValue $ws "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D50" "0.0524"
Set-TextValue $ws "E50" "  -0.05%  "
Set-TextValue $ws "D51" "39.72"
Set-TextValue $ws "E51" "  +17.19%  "
